# The source data added one more weekly record (new commit: "Fruta / hortaliza, semanal").
# A brand-new row is inserted right before the existing row 51, shifting every
# subsequent record down by one row (old row 51 -> new row 52, ... old row 153 -> new row 154).
# The new row carries the same product/market/price-unit profile as the record
# that ends up just below it (old row 51, i.e. the "3 kilos / Región Metropolitana"
# pattern), but with a new observation date (44469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51; Excel shifts rows 51:153 down to 52:154
# and the worksheet dimension grows to A1:R154 automatically.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new record.
$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 44469
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112044
$ws.Range("G51").Value = "Perejil"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 140
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = 4500
$ws.Range("N51").Value = "$/docena de atados (3 kilos)"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 1500
$ws.Range("Q51").Value = 3
$ws.Range("R51").Value = "Hortaliza"
